$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 97, shifting existing rows 97:234 down to 98:235
$ws.Rows("97:97").Insert(1)

# Populate the newly inserted row 97 with the new data record
$ws.Range("A97").Value = 10
$ws.Range("B97").Value = 'Vega Modelo de Temuco'
$ws.Range("C97").Value = 'La Araucanía'
$ws.Range("D97").Value = 44803
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = 100112005
$ws.Range("G97").Value = 'Puerro'
$ws.Range("H97").Value = 'Azul de Maquehue'
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 30
$ws.Range("K97").Value = 16000
$ws.Range("L97").Value = 16000
$ws.Range("M97").Value = 16000
$ws.Range("N97").Value = '$/docena de paquetes'
$ws.Range("O97").Value = 'Provincia de Cautín'
$ws.Range("P97").Value = 1333
$ws.Range("Q97").Value = 12
$ws.Range("R97").Value = 'Hortaliza'
